# Apply the "added architecture and updated flow diagram" edit to slide 1.
# (Note: the presentation also carries a cosmetic, auto-updating
# "datetimeFigureOut" field -- the date PowerPoint stamps into the
# master/layout footers whenever the file is resaved. It is not part of
# the actual diagram edit and the COM TextRange.Text setter would convert
# that <a:fld> into a plain run if touched, so it is intentionally left
# alone here.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Picture 8 (MIMIC graph icon, top) : shift right ---
$pic8 = $s.Shapes.Item(2)
$pic8.Left = 3066177 / 12700

# --- Picture 10 (MIMIC graph icon, bottom) : shift right ---
$pic10 = $s.Shapes.Item(3)
$pic10.Left = 3075366 / 12700

# --- TextBox 11 ("Patient Graph with Diagnosis, etc.") : shift right ---
$tb11 = $s.Shapes.Item(4)
$tb11.Left = 2618341 / 12700

# --- TextBox 14 ("Patient Similarity Graph" -> "Patient Similarity Pairings") : shift right + rename ---
$tb14 = $s.Shapes.Item(5)
$tb14.Left = 2663090 / 12700
$tb14.TextFrame.TextRange.Text = "Patient Similarity Pairings"

# --- TextBox 25 ("Graph Neural Network" -> "Relational Graph Convolutional Network") : resize + rename ---
$tb25 = $s.Shapes.Item(9)
$tb25.Width = 1351295 / 12700
$tb25.Height = 646331 / 12700
$tb25.TextFrame.TextRange.Text = "Relational Graph Convolutional Network"
